# "Update Run as FireFox" -- add new quiz questions (Question 7 and Question 38)
# to the "Exam Sample A" and "Exam Sample B" worksheets.
#
# Helper: write a truly-empty TEXT cell (matches the workbook's existing
# convention of using the shared empty string "" rather than leaving the
# cell completely untouched). A bare Value = "" collapses to a blank
# Number cell in this engine, so we go through the classic "force text"
# apostrophe prefix and then reset the style back to Normal so no stray
# cell-level style index is left behind.
function Set-EmptyTextCell($cell) {
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# Helper: make sure a blank cell is materialised in the sheet XML (an
# empty <c r=".."/> tag) without altering its style. Touching Font.Bold
# with its own default value (False) is enough to "dirty" the cell while
# being a complete no-op, so no new style record is created.
function Set-TouchedBlankCell($cell) {
    $cell.Font.Bold = $false
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Exam Sample A": append Question 7 (rows 43-52)
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Exam Sample A")

# Row 43: blank separator row (three empty-string cells)
Set-EmptyTextCell $wsA.Cells.Item(43, 1)
Set-EmptyTextCell $wsA.Cells.Item(43, 2)
Set-EmptyTextCell $wsA.Cells.Item(43, 3)

# Row 44: header row ("Questions" / "Answer"), first cell stays blank
Set-TouchedBlankCell $wsA.Cells.Item(44, 1)
$wsA.Cells.Item(44, 2).Value = "Questions"
$wsA.Cells.Item(44, 3).Value = "Answer"

$q7Label = "Question 7:"
$q7Text = "Question #7`nWhich of the following activities is part of the main activity ""test analysis"" in the test process?"
$q7AnswerA = "a) Identifying any required infrastructure and tools."
$q7AnswerB = "b) Creating test suites from test scripts."
$q7AnswerC = "c) Analyzing lessons learned for process improvement."
$q7AnswerD = "d) Evaluating the test basis for testability."

$q7Rows = @(
    @(45, $q7AnswerA),
    @(46, $q7AnswerB),
    @(47, $q7AnswerC),
    @(48, $q7AnswerD),
    @(49, $q7AnswerA),
    @(50, $q7AnswerB),
    @(51, $q7AnswerC),
    @(52, $q7AnswerD)
)

foreach ($row in $q7Rows) {
    $r = $row[0]
    $wsA.Cells.Item($r, 1).Value = $q7Label
    $wsA.Cells.Item($r, 2).Value = $q7Text
    $wsA.Cells.Item($r, 3).Value = $row[1]
}

# ---------------------------------------------------------------------
# Sheet "Exam Sample B": append Question 38 (rows 11-15)
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Exam Sample B")

# Row 11: blank separator row (three empty-string cells)
Set-EmptyTextCell $wsB.Cells.Item(11, 1)
Set-EmptyTextCell $wsB.Cells.Item(11, 2)
Set-EmptyTextCell $wsB.Cells.Item(11, 3)

# Row 12: header row ("Questions" / "Answer"), first cell stays blank
Set-TouchedBlankCell $wsB.Cells.Item(12, 1)
$wsB.Cells.Item(12, 2).Value = "Questions"
$wsB.Cells.Item(12, 3).Value = "Answer"

$q38Label = "Question 38:"
$q38Text = "Question #38`nYou are performing system testing of a train reservation system. Based on the test cases performed, you have noticed that the system occasionally reports that no trains are available, although this should actually be the case. You have provided the developers with a summary of the defect and the version of the tested system. They recognize the urgency of the defect and are now waiting for you to provide further details.`nIn addition to the information already provided, the following additional information is given:`n1. Degree of impact (severity) of the defect.`n2. Identification of the test item.`n3. Details of the test environment.`n4. Urgency/priority to fix.`n5. Actual results.`n6. Reference to test case specification.`nWhich of this information is most useful to include in the defect report?"

$wsB.Cells.Item(13, 1).Value = $q38Label
$wsB.Cells.Item(13, 2).Value = $q38Text
$wsB.Cells.Item(13, 3).Value = "1. Degree of impact (severity) of the defect."

$wsB.Cells.Item(14, 1).Value = $q38Label
$wsB.Cells.Item(14, 2).Value = $q38Text
$wsB.Cells.Item(14, 3).Value = "a) 1, 2, 6"

$wsB.Cells.Item(15, 1).Value = $q38Label
$wsB.Cells.Item(15, 2).Value = $q38Text
$wsB.Cells.Item(15, 3).Value = "b) 1, 4, 5, 6"
